# Generate Report for Handoff
#
# Replaces the two "handed back" localization file entries
# (11997f48-...md / af1f09fe-...md) with a fresh pair of files that are
# "Ready for handoff" (026cd78d-...md / ffffec5bbaeb-...md), refreshes the
# handoff file hashes/timestamps, drops the now-stale "Latest Target
# File"/"Latest Handback File" columns (E/F) on the per-locale sheets, and
# updates every hyperlink's display text + target to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New identifiers / values used throughout the workbook.
# ---------------------------------------------------------------------
$oldFile1 = "11997f48-c54a-44ad-a860-376efe8dc576.md"
$oldFile2 = "af1f09fe-b970-4b01-9e3a-5bd87979b265.md"

$newFile1 = "026cd78d-c75f-434c-9286-2379c3ec760f.md"
$newFile2 = "ffffec5bbaeb-f49e-47f5-a051-5415bf9173e7.md"

$newStatus = "Ready for handoff"

$newHandoffZhCn = "026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.zh-cn.xlf"
$newHandoffDeDe = "026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.de-de.xlf"

$newHandoffDtZhCn = "2016-03-08 06:34:17"
$newHandoffDtDeDe = "2016-03-08 06:34:20"

$emptyDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: filenames + status text, and hyperlink display text.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newFile1
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus

$ov.Range("A3").Value = $newFile2
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/026cd78d-c75f-434c-9286-2379c3ec760f.md", "", "", $newFile1)
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/ffffec5bbaeb-f49e-47f5-a051-5415bf9173e7.md", "", "", $newFile2)
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newFile1
$zh.Range("B2").Value = $newStatus
$zh.Range("C2").Value = $newHandoffZhCn
$zh.Range("D2").Value = $newHandoffDtZhCn
$zh.Range("E2").Clear()
$zh.Range("F2").Clear()
$zh.Range("G2").Value = $emptyDate
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = $newFile2
$zh.Range("B3").Value = $newStatus
$zh.Range("C3").Value = $newHandoffZhCn
$zh.Range("D3").Value = $newHandoffDtZhCn
$zh.Range("E3").Clear()
$zh.Range("F3").Clear()
$zh.Range("G3").Value = $emptyDate
$zh.Range("H3").Value = "Include"

$zh.Range("D4").Value = $emptyDate
$zh.Range("G4").Value = $emptyDate
$zh.Range("H4").Value = "Ignored"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/026cd78d-c75f-434c-9286-2379c3ec760f.md", "", "", $newFile1)
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/781e63214eebd32ee9937c4860216e93467692d6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newHandoffZhCn", "", "", $newHandoffZhCn)
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/ffffec5bbaeb-f49e-47f5-a051-5415bf9173e7.md", "", "", $newFile2)
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/781e63214eebd32ee9937c4860216e93467692d6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newHandoffZhCn", "", "", $newHandoffZhCn)
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newFile1
$de.Range("B2").Value = $newStatus
$de.Range("C2").Value = $newHandoffDeDe
$de.Range("D2").Value = $newHandoffDtDeDe
$de.Range("E2").Clear()
$de.Range("F2").Clear()
$de.Range("G2").Value = $emptyDate
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = $newFile2
$de.Range("B3").Value = $newStatus
$de.Range("C3").Value = $newHandoffDeDe
$de.Range("D3").Value = $newHandoffDtDeDe
$de.Range("E3").Clear()
$de.Range("F3").Clear()
$de.Range("G3").Value = $emptyDate
$de.Range("H3").Value = "Include"

$de.Range("D4").Value = $emptyDate
$de.Range("G4").Value = $emptyDate
$de.Range("H4").Value = "Ignored"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/026cd78d-c75f-434c-9286-2379c3ec760f.md", "", "", $newFile1)
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90892827be9efb607d9daaf1c387f29b8ef3c264/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newHandoffDeDe", "", "", $newHandoffDeDe)
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/e2e/ffffec5bbaeb-f49e-47f5-a051-5415bf9173e7.md", "", "", $newFile2)
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90892827be9efb607d9daaf1c387f29b8ef3c264/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newHandoffDeDe", "", "", $newHandoffDeDe)
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2dd4f6df1475f5d8c53f4e29b2110b8eef71fe18/.localization-config", "", "", ".localization-config")
